# Actualización automática 2025-07-10 14:00:09
$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("H2").Value = 811.8
$ws1.Range("H32").Value = "1 de 30"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F2").Value = 811.8
$ws2.Range("F32").Value = 1048.99

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D6").Value = 811.8
$ws3.Range("E6").Value = 788.2
$ws3.Range("F6").Value = 0.507375

$ws3.Range("D18").Value = 1038.71
$ws3.Range("E18").Value = 32896.00607548726
$ws3.Range("F18").Value = 0.03060906705950937
